# Refine metadata into its own "metadata" worksheet, and bump the
# "time_taken" query timestamps on the "data" sheet to reflect the
# re-run that produced the metadata tab.

$wb = $excel.ActiveWorkbook

# --- 1. Update the per-gene "time_taken" timestamps on the data sheet ---
$dataWs = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:21:35.013852",
    "2021-10-05 14:21:35.013860",
    "2021-10-05 14:21:35.013863",
    "2021-10-05 14:21:35.013865",
    "2021-10-05 14:21:35.013868",
    "2021-10-05 14:21:35.013871",
    "2021-10-05 14:21:35.013873",
    "2021-10-05 14:21:35.013876",
    "2021-10-05 14:21:35.013878",
    "2021-10-05 14:21:35.013881",
    "2021-10-05 14:21:35.013883",
    "2021-10-05 14:21:35.013886",
    "2021-10-05 14:21:35.013888",
    "2021-10-05 14:21:35.013891",
    "2021-10-05 14:21:35.013893",
    "2021-10-05 14:21:35.013896",
    "2021-10-05 14:21:35.013898",
    "2021-10-05 14:21:35.013901",
    "2021-10-05 14:21:35.013903",
    "2021-10-05 14:21:35.013906",
    "2021-10-05 14:21:35.013909",
    "2021-10-05 14:21:35.013911",
    "2021-10-05 14:21:35.013914",
    "2021-10-05 14:21:35.013916",
    "2021-10-05 14:21:35.013919",
    "2021-10-05 14:21:35.013921",
    "2021-10-05 14:21:35.013924",
    "2021-10-05 14:21:35.013926"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataWs.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- 2. Add the new "metadata" worksheet, placed right after "data" ---
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Header row - reuse the same bordered/bold/centered header style as
# the "data" sheet (columns B:F carry it there already).
$dataWs.Range("B1:F1").Copy($metaWs.Range("B1:F1"))
$dataWs.Range("B1").Copy($metaWs.Range("G1"))

$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

# Row 2 - the panel-level metadata record. A2 reuses the index-column
# style, the rest are plain cells (matching the "data" sheet's layout).
$dataWs.Range("A2").Copy($metaWs.Range("A2"))
$metaWs.Range("A2").Value = 0

$metaWs.Range("B2").Value = "Mitochondrial DNA maintenance disorder"
$metaWs.Range("C2").Value = 533

# data_version must stay textual ("1.4", not the float 1.4) - force
# text storage, then drop the format so the cell ends up unstyled.
$metaWs.Range("D2").NumberFormat = "@"
$metaWs.Range("D2").Value = "1.4"
$metaWs.Range("D2").ClearFormats()

$metaWs.Range("E2").Value = "2021-05-19T12:55:34.455564Z"
$metaWs.Range("F2").Value = "2021-10-05 14:21:35.010229"
$metaWs.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/533/?format=json"
